$d = $word.ActiveDocument

# Use Track Changes so the engine does not silently coalesce the
# surrounding (untouched) runs into the run we are editing; we accept
# all the resulting tracked insertions at the end so the final document
# contains plain (non-tracked) runs.
$d.TrackRevisions = $true

# --- Change 1 -----------------------------------------------------
# "I would like for the library to allow a large degree of customization while be"
# becomes three runs:
#   "I would like for the library to allow a large degree of customization "
#   "(by programmers; in other words, to be an extensible API) "
#   "while be"
$rng1 = $d.Content
$rng1.Find.Execute("while be", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1.Collapse(1)
$rng1.InsertBefore("(by programmers; in other words, to be an extensible API) ")

# --- Change 2 -----------------------------------------------------
# "his library is for my use, ... Each subsystem is independent, and it's
#  possible to use components individually. "
# becomes three runs:
#   "his library is for my use, ... Each subsystem is independent"
#   " (except that they all depend on utility)"
#   ", and it's possible to use components individually. "
$rng2 = $d.Content
$rng2.Find.Execute("Each subsystem is independent", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$rng2.InsertAfter(" (except that they all depend on utility)")

# Accept the tracked insertions so the saved document has ordinary runs.
$d.TrackRevisions = $false
$d.AcceptAllRevisions()
